$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text, matching the workbook's existing
# inline-string cells (many values like "1.001" or "5.182" would otherwise be
# auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.349.70'
$ws.Range("E2").Value = '  -1.20%  '

$ws.Range("D3").Value = '1.891.90'
$ws.Range("E3").Value = '  -1.22%  '

$ws.Range("D5").Value = '238.31'
$ws.Range("E5").Value = '  -1.19%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.17%  '

$ws.Range("E8").Value = '  -3.38%  '

$ws.Range("D9").Value = '0.06606'
$ws.Range("E9").Value = '  -2.41%  '

$ws.Range("D10").Value = '1.902.51'
$ws.Range("E10").Value = '  -0.69%  '

$ws.Range("D11").Value = '16.93'
$ws.Range("E11").Value = '  -1.55%  '

$ws.Range("D12").Value = '0.07386'
$ws.Range("E12").Value = '  +0.77%  '

$ws.Range("D13").Value = '5.182'
$ws.Range("E13").Value = '  -0.72%  '

$ws.Range("D14").Value = '87.66'
$ws.Range("E14").Value = '  -1.08%  '

$ws.Range("D15").Value = '0.6637'
$ws.Range("E15").Value = '  -1.85%  '

$ws.Range("D16").Value = '30.340.87'
$ws.Range("E16").Value = '  -1.25%  '

$ws.Range("D17").Value = '13.44'
$ws.Range("E17").Value = '  -0.97%  '

$ws.Range("D18").Value = '0.000007773'
$ws.Range("E18").Value = '  -2.61%  '

$ws.Range("D19").Value = '0.9999'
$ws.Range("E19").Value = '  -0.12%  '

$ws.Range("D20").Value = '5.501'
$ws.Range("E20").Value = '  +1.15%  '

$ws.Range("D21").Value = '2.146.86'
$ws.Range("E21").Value = '  -0.84%  '

$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.16%  '

$ws.Range("D23").Value = '192.11'
$ws.Range("E23").Value = '  -3.28%  '

$ws.Range("D24").Value = '6.198'
$ws.Range("E24").Value = '  -2.51%  '

$ws.Range("D25").Value = '9.457'
$ws.Range("E25").Value = '  -2.31%  '

$ws.Range("D26").Value = '165.43'
$ws.Range("E26").Value = '  +2.13%  '

$ws.Range("D27").Value = '18.21'
$ws.Range("E27").Value = '  -2.56%  '

$ws.Range("D28").Value = '1.963'
$ws.Range("E28").Value = '  -0.12%  '

$ws.Range("D29").Value = '1.450'
$ws.Range("E29").Value = '  -0.59%  '

$ws.Range("D30").Value = '4.265'
$ws.Range("E30").Value = '  -2.13%  '

$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("D32").Value = '4.050'
$ws.Range("E32").Value = '  -0.63%  '

$ws.Range("D33").Value = '0.05089'
$ws.Range("E33").Value = '  -3.94%  '

$ws.Range("D34").Value = '0.7352'
$ws.Range("E34").Value = '  -1.38%  '

$ws.Range("D35").Value = '1.145'
$ws.Range("E35").Value = '  +1.62%  '

$ws.Range("D36").Value = '2.711'
$ws.Range("E36").Value = '  -0.15%  '

$ws.Range("D37").Value = '0.01825'
$ws.Range("E37").Value = '  -1.94%  '

$ws.Range("D38").Value = '2.647'
$ws.Range("E38").Value = '  -3.04%  '

$ws.Range("D39").Value = '0.9206'
$ws.Range("E39").Value = '  -1.09%  '

$ws.Range("D40").Value = '2.084'
$ws.Range("E40").Value = '  -0.25%  '

$ws.Range("D41").Value = '106.72'
$ws.Range("E41").Value = '  -0.67%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '5.911'
$ws.Range("E42").Value = '  -0.78%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4345'
$ws.Range("E43").Value = '  -3.73%  '

$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '0.1369'
$ws.Range("E45").Value = '  -2.17%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.630'
$ws.Range("E46").Value = '  -1.28%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '1.570'
$ws.Range("E47").Value = '  +8.05%  '

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '65.52'
$ws.Range("E48").Value = '  -9.12%  '

$ws.Range("D49").Value = '9.033'
$ws.Range("E49").Value = '  -0.93%  '

$ws.Range("D50").Value = '34.28'
$ws.Range("E50").Value = '  -3.84%  '

$ws.Range("D51").Value = '0.05782'
$ws.Range("E51").Value = '  -2.01%  '
